# Apply updated cryptocurrency price/volume data (D and E columns) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '53.560.58'
Set-TextValue 'E2' '  -11.51%  '
Set-TextValue 'D3' '2.324.37'
Set-TextValue 'E3' '  -19.78%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '440.68'
Set-TextValue 'E5' '  -16.41%  '
Set-TextValue 'D6' '124.09'
Set-TextValue 'E6' '  -13.25%  '
Set-TextValue 'E7' '  -0.19%  '
Set-TextValue 'D8' '0.477'
Set-TextValue 'E8' '  -14.20%  '
Set-TextValue 'D9' '2.333.20'
Set-TextValue 'E9' '  -19.72%  '
Set-TextValue 'D10' '5.33'
Set-TextValue 'E10' '  -11.63%  '
Set-TextValue 'D11' '0.0919'
Set-TextValue 'E11' '  -14.70%  '
Set-TextValue 'D12' '0.309'
Set-TextValue 'E12' '  -14.57%  '
Set-TextValue 'E13' '  -3.32%  '
Set-TextValue 'D14' '2.728.26'
Set-TextValue 'E14' '  -19.73%  '
Set-TextValue 'D15' '53.532.27'
Set-TextValue 'E15' '  -11.55%  '
Set-TextValue 'D16' '18.92'
Set-TextValue 'E16' '  -16.90%  '
Set-TextValue 'E17' '  -14.17%  '
Set-TextValue 'D18' '2.348.88'
Set-TextValue 'E18' '  -19.16%  '
Set-TextValue 'D19' '3.96'
Set-TextValue 'E19' '  -21.33%  '
Set-TextValue 'D20' '299.20'
Set-TextValue 'E20' '  -17.16%  '
Set-TextValue 'D21' '9.21'
Set-TextValue 'E21' '  -21.31%  '
Set-TextValue 'D22' '0.998'
Set-TextValue 'E22' '  -0.18%  '
Set-TextValue 'D23' '5.61'
Set-TextValue 'E23' '  -1.17%  '
Set-TextValue 'D24' '5.40'
Set-TextValue 'E24' '  -18.81%  '
Set-TextValue 'D25' '55.58'
Set-TextValue 'E25' '  -14.09%  '
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  -0.14%  '
Set-TextValue 'D27' '0.151'
Set-TextValue 'E27' '  -16.17%  '
Set-TextValue 'D28' '0.368'
Set-TextValue 'E28' '  -19.04%  '
Set-TextValue 'D29' '6.93'
Set-TextValue 'E29' '  -12.03%  '
Set-TextValue 'D30' '0.997'
Set-TextValue 'E30' '  -0.22%  '
Set-TextValue 'E31' '  -17.44%  '
Set-TextValue 'D32' '145.35'
Set-TextValue 'E32' '  -4.89%  '
Set-TextValue 'D33' '17.21'
Set-TextValue 'E33' '  -12.74%  '
Set-TextValue 'E34' '  -19.76%  '
Set-TextValue 'D35' '4.68'
Set-TextValue 'E35' '  -16.12%  '
Set-TextValue 'D36' '3.55'
Set-TextValue 'E36' '  -18.84%  '
Set-TextValue 'D37' '0.824'
Set-TextValue 'E37' '  -18.21%  '
Set-TextValue 'E38' '  -16.41%  '
Set-TextValue 'D39' '33.27'
Set-TextValue 'E39' '  -11.62%  '
Set-TextValue 'D40' '0.995'
Set-TextValue 'E40' '  -0.29%  '
Set-TextValue 'E41' '  -0.48%  '
Set-TextValue 'D42' '1.943.47'
Set-TextValue 'E42' '  -15.06%  '
Set-TextValue 'D43' '3.14'
Set-TextValue 'E43' '  -15.62%  '
Set-TextValue 'E44' '  -18.06%  '
Set-TextValue 'D45' '0.0495'
Set-TextValue 'E45' '  -14.76%  '
Set-TextValue 'D46' '0.521'
Set-TextValue 'E46' '  -19.45%  '
Set-TextValue 'D47' '0.0209'
Set-TextValue 'E47' '  -11.85%  '
Set-TextValue 'D48' '0.0831'
Set-TextValue 'E48' '  -10.05%  '
Set-TextValue 'D49' '15.93'
Set-TextValue 'E49' '  -22.03%  '
Set-TextValue 'D50' '3.96'
Set-TextValue 'E50' '  -20.72%  '
Set-TextValue 'D51' '4.64'
Set-TextValue 'E51' '  -3.89%  '
